# Week 23 Timesheet - Caroline Chang
# "panels updated, inspect is a m mess."
#
# Sunday (column H) work was logged for the "Sponsor Work" task (row 13):
# 3 hours on Sunday, which also feeds the "Daily Total" row (row 14) and
# bumps both rows' weekly "Daily Total" column (I) from 3.5 to 6.5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Sunday entry for "Sponsor Work" (row 13).
$ws.Range("H13").Value = 3

# Daily total for "Sponsor Work" row now includes the new Sunday hours.
$ws.Range("I13").Value = 6.5

# The "Daily Total" row (14) picks up the same Sunday total...
$ws.Range("H14").Value = 3

# ...and its own weekly total grows to match.
$ws.Range("I14").Value = 6.5

# Reflect the author's final on-sheet selection (K10:K11, anchored at K11).
$ws.Range("K10:K11").Select()
